# update utils and result
# Adds a new "LRU" column (J) with timing values alongside the existing
# A/B/C (Trace/LRU/Random) columns, and moves the active selection to J31.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header: same label as B1 ("LRU")
$ws.Range("J1").Value = "LRU"

$jValues = @(
    11.305999999999999,
    88.179000000000002,
    98.787999999999997,
    88.962000000000003,
    82.284999999999997,
    111.264,
    95.216999999999999,
    86.534999999999997,
    91.88,
    93.781999999999996,
    89.251999999999995,
    75.444999999999993,
    86.007000000000005,
    100.492,
    56.247,
    102.956,
    96.728999999999999,
    92.484999999999999,
    97.311000000000007,
    89.581000000000003,
    90.14,
    86.165000000000006,
    80.817999999999998,
    98.933000000000007,
    91.725999999999999,
    96.027000000000001,
    89.194000000000003,
    94.662999999999997,
    52.363999999999997
)

for ($i = 0; $i -lt $jValues.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 10).Value = $jValues[$i]
}

# Reflect the author's final cursor position/selection
$ws.Range("J31").Select()
